# "Generate Report for Archive"
# The localization status report moved from "Ready for handoff" to
# "In Translation" for both locales, and the now-narrower status text let
# the status columns shrink (handoff-date columns on Overview, and the
# Status column on each locale sheet).

$wb = $excel.ActiveWorkbook

# --- Update the status text on every sheet that reports it -----------------
$overview = $wb.Sheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Sheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- Shrink the status columns to match the shorter text --------------------
$overview.Columns("E:F").ColumnWidth = 12.5
$zhcn.Columns("C:C").ColumnWidth = 12.5
$dede.Columns("C:C").ColumnWidth = 12.5
